$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '71.028.35'
$ws.Range('E2').Value = '  +2.57%  '
$ws.Range('D3').Value = '3.807.07'
$ws.Range('E3').Value = '  +0.87%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = '699.27'
$ws.Range('E5').Value = '  +10.68%  '
$ws.Range('D6').Value = '172.41'
$ws.Range('E6').Value = '  +3.32%  '
$ws.Range('D7').Value = '3.806.92'
$ws.Range('E7').Value = '  +0.93%  '
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('E9').Value = '  +0.77%  '
$ws.Range('E10').Value = '  +3.01%  '
$ws.Range('D11').Value = '7.53'
$ws.Range('E11').Value = '  +11.85%  '
$ws.Range('E12').Value = '  +0.54%  '
$ws.Range('E13').Value = '  +9.80%  '
$ws.Range('D14').Value = '36.27'
$ws.Range('E14').Value = '  +3.45%  '
$ws.Range('D15').Value = '4.451.07'
$ws.Range('E15').Value = '  +0.92%  '
$ws.Range('D16').Value = '3.798.80'
$ws.Range('E16').Value = '  +0.97%  '
$ws.Range('D17').Value = '70.976.16'
$ws.Range('E17').Value = '  +2.53%  '
$ws.Range('D18').Value = '17.83'
$ws.Range('E18').Value = '  +1.34%  '
$ws.Range('D19').Value = '7.21'
$ws.Range('E19').Value = '  +2.79%  '
$ws.Range('E20').Value = '  +1.14%  '
$ws.Range('D21').Value = '11.24'
$ws.Range('E21').Value = '  +18.05%  '
$ws.Range('D22').Value = '480.26'
$ws.Range('E22').Value = '  +3.68%  '
$ws.Range('D23').Value = '0.716'
$ws.Range('E23').Value = '  +1.35%  '
$ws.Range('D24').Value = '83.97'
$ws.Range('E24').Value = '  +1.55%  '
$ws.Range('E25').Value = '  +0.86%  '
$ws.Range('D26').Value = '12.37'
$ws.Range('E26').Value = '  +2.41%  '
$ws.Range('E27').Value = '  +1.20%  '
$ws.Range('E28').Value = '  +3.95%  '
$ws.Range('D29').Value = '3.959.20'
$ws.Range('E29').Value = '  +0.92%  '
$ws.Range('E30').Value = '  -0.05%  '
$ws.Range('E31').Value = '  +14.97%  '
$ws.Range('E32').Value = '  -0.53%  '
$ws.Range('E33').Value = '  +6.69%  '
$ws.Range('D34').Value = '29.55'
$ws.Range('E34').Value = '  +3.84%  '
$ws.Range('E35').Value = '  +5.75%  '
$ws.Range('D36').Value = '9.22'
$ws.Range('E36').Value = '  +3.17%  '
$ws.Range('D37').Value = '0.999'
$ws.Range('E37').Value = '  +0.03%  '
$ws.Range('D38').Value = '3.757.92'
$ws.Range('E39').Value = '  +1.27%  '
$ws.Range('D40').Value = '3.48'
$ws.Range('E40').Value = '  +5.86%  '
$ws.Range('D41').Value = '5.98'
$ws.Range('E41').Value = '  +3.08%  '
$ws.Range('E42').Value = '  +11.58%  '
$ws.Range('D43').Value = '0.000327'
$ws.Range('E43').Value = '  +22.34%  '
$ws.Range('D44').Value = '0.969'
$ws.Range('E44').Value = '  +0.75%  '
$ws.Range('E45').Value = '  -0.02%  '
$ws.Range('E46').Value = '  +0.00%  '
$ws.Range('D47').Value = '45.51'
$ws.Range('E47').Value = '  +5.17%  '
$ws.Range('D48').Value = '160.63'
$ws.Range('E48').Value = '  +1.75%  '
$ws.Range('D49').Value = '49.30'
$ws.Range('E49').Value = '  +5.41%  '
$ws.Range('E50').Value = '  -1.31%  '
$ws.Range('E51').Value = '  +1.22%  '
